$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2: add the google-phone "treated on" timestamp (text, matches existing F-column text format)
$ws.Range("F2").Value = "29/01/2025 18:42"

# Row 2 grows taller once the wrapped text lands in F2 (matches Excel's own re-layout)
$ws.Rows(2).RowHeight = 28.5
